$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1. Move the existing text box ("文本框 2") to its new position.
$existing = $s.Shapes.Item(1)
$existing.Left = 2096527 / 12700
$existing.Top = 4916556 / 12700

# 2. Add the new text box ("文本框 1") with the explanatory paragraphs.
$newBox = $s.Shapes.AddTextbox(1, 2114159 / 12700, 5579164 / 12700, 7491153 / 12700, 923330 / 12700)
$newBox.Name = "文本框 1"
$newBox.TextFrame.WordWrap = $false

$tr = $newBox.TextFrame.TextRange
$tr.Text = "中间一列是路径参数，eg:localhost:8088/getStu/1203，1203是参数id，`r但是在这里表示为路径的一部分。`r通常是用：localost:8088/getStu?id=1203来请求的，可读性也更强"
